$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("quality_comparison")
$ws2 = $wb.Worksheets.Item("computational_comparison")

# --- quality_comparison sheet ---
# C1: give it a top+bottom border only (drop the inherited bold/box style).
# Order matters: apply Bottom/Right before Top so no partial state along the
# way happens to coincide with one of the pre-existing border combos (which
# would otherwise leave a spurious, unused cellXf behind).
$c = $ws1.Range("C1")
$c.ClearFormats()
$c.Borders(9).LineStyle = 1
$c.Borders(8).LineStyle = 1

# D1: top+bottom+right border only
$d = $ws1.Range("D1")
$d.ClearFormats()
$d.Borders(9).LineStyle = 1
$d.Borders(10).LineStyle = 1
$d.Borders(8).LineStyle = 1

# C2 header text: anonymize "fedcore" -> "approach"
$ws1.Range("C2").Value = "approach"

# --- computational_comparison sheet ---
# C1: top+bottom border only
$c2 = $ws2.Range("C1")
$c2.ClearFormats()
$c2.Borders(9).LineStyle = 1
$c2.Borders(8).LineStyle = 1

# D1: top+bottom+right border only
$d2 = $ws2.Range("D1")
$d2.ClearFormats()
$d2.Borders(9).LineStyle = 1
$d2.Borders(10).LineStyle = 1
$d2.Borders(8).LineStyle = 1

# F1: top+bottom border only
$f2 = $ws2.Range("F1")
$f2.ClearFormats()
$f2.Borders(9).LineStyle = 1
$f2.Borders(8).LineStyle = 1

# G1: top+bottom+right border only
$g2 = $ws2.Range("G1")
$g2.ClearFormats()
$g2.Borders(9).LineStyle = 1
$g2.Borders(10).LineStyle = 1
$g2.Borders(8).LineStyle = 1

# C2 / F2 header text: anonymize "fedcore" -> "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# G5 was a stray empty cell; remove it entirely
$ws2.Range("G5").ClearContents()
